$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the formatting from the last date column (BE) onto the new
# "30-ago" column (BF) before filling in its header + values.
$ws.Range("BE1:BE11").Copy()
$ws.Range("BF1:BF11").PasteSpecial(-4122)

$ws.Range("BF1").Value = "30-ago"

$values = @(14, 10, 10, 13, 9, 16, 12, 17, 18, 9)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 58).Value = $values[$i]
}

# Update the active cell selection to match the saved view
$ws.Range("BG8").Select()
